$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3849.95
$ws.Range("I113").Value = 3418.182
$ws.Range("J113").Value = 4377.6665
$ws.Range("K113").Value = 3418.182
$ws.Range("L113").Value = 4377.6665
$ws.Range("M113").Value = -164.1819999999998
$ws.Range("N113").Value = -10885.6665
$ws.Range("H116").Value = 4923
$ws.Range("I116").Value = 2795
$ws.Range("J116").Value = 5189
$ws.Range("K116").Value = 2795
$ws.Range("L116").Value = 5189
$ws.Range("M116").Value = 647
$ws.Range("N116").Value = -12073

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1152.3939
$ws.Range("I2").Value = 1175.0952
$ws.Range("J2").Value = 1112.6666
$ws.Range("K2").Value = 1175.0952
$ws.Range("L2").Value = 1112.6666
$ws.Range("M2").Value = -1062.0952
$ws.Range("N2").Value = -1338.6666
$ws.Range("H32").Value = 19630.035
$ws.Range("I32").Value = 20338.74
$ws.Range("K32").Value = 20338.74
$ws.Range("M32").Value = -20051.74
$ws.Range("H45").Value = 1822.7333
$ws.Range("I45").Value = 1970.25
$ws.Range("J45").Value = 1654.1428
$ws.Range("K45").Value = 1970.25
$ws.Range("L45").Value = 1654.1428
$ws.Range("M45").Value = -1593.25
$ws.Range("N45").Value = -2408.1428
$ws.Range("H61").Value = 2155.8
$ws.Range("I61").Value = 1712.1724
$ws.Range("J61").Value = 4300
$ws.Range("K61").Value = 1712.1724
$ws.Range("L61").Value = 4300
$ws.Range("M61").Value = -1500.1724
$ws.Range("N61").Value = -4724
$ws.Range("H116").Value = 1152.3939
$ws.Range("I116").Value = 1175.0952
$ws.Range("J116").Value = 1112.6666
$ws.Range("K116").Value = 1175.0952
$ws.Range("L116").Value = 1112.6666
$ws.Range("M116").Value = 1118.9048
$ws.Range("N116").Value = -5700.6666
$ws.Range("H132").Value = 11408.862
$ws.Range("I132").Value = 1410.9231
$ws.Range("J132").Value = 43902.168
$ws.Range("K132").Value = 4232.7693
$ws.Range("L132").Value = 131706.504
$ws.Range("M132").Value = -1702.7693
$ws.Range("N132").Value = -136766.504
$ws.Range("H136").Value = 2155.8
$ws.Range("I136").Value = 1712.1724
$ws.Range("J136").Value = 4300
$ws.Range("K136").Value = 5136.5172
$ws.Range("L136").Value = 12900
$ws.Range("M136").Value = -2586.5172
$ws.Range("N136").Value = -18000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1152.3939
$ws.Range("I3").Value = 1175.0952
$ws.Range("J3").Value = 1112.6666
$ws.Range("K3").Value = 1175.0952
$ws.Range("L3").Value = 1112.6666
$ws.Range("M3").Value = -1061.0952
$ws.Range("N3").Value = -1340.6666
$ws.Range("H99").Value = 1467.4375
$ws.Range("I99").Value = 937.9
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 937.9
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = 560.1
$ws.Range("N99").Value = -5346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12245.066
$ws.Range("I31").Value = 16022.6
$ws.Range("K31").Value = 16022.6
$ws.Range("M31").Value = -15727.6
$ws.Range("H34").Value = 12245.066
$ws.Range("I34").Value = 16022.6
$ws.Range("K34").Value = 16022.6
$ws.Range("M34").Value = -15820.6
$ws.Range("H99").Value = 14709804
$ws.Range("I99").Value = 3184.1738
$ws.Range("J99").Value = 45460010
$ws.Range("K99").Value = 3184.1738
$ws.Range("L99").Value = 45460010
$ws.Range("M99").Value = -1686.1738
$ws.Range("N99").Value = -45463006
$ws.Range("H126").Value = 14709804
$ws.Range("I126").Value = 3184.1738
$ws.Range("J126").Value = 45460010
$ws.Range("K126").Value = 9552.5214
$ws.Range("L126").Value = 136380030
$ws.Range("M126").Value = -7082.5214
$ws.Range("N126").Value = -136384970
$ws.Range("H127").Value = 39926.668
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 39926.668
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 39926.668
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -49846.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 302
$ws.Range("I14").Value = 302
$ws.Range("K14").Value = 906
$ws.Range("M14").Value = -733
$ws.Range("H118").Value = 100002330
$ws.Range("J118").Value = 3833.3333
$ws.Range("L118").Value = 11499.9999
$ws.Range("N118").Value = -13985.9999
$ws.Range("H131").Value = 730.46
$ws.Range("J131").Value = 747.8172
$ws.Range("L131").Value = 2243.4516
$ws.Range("N131").Value = -12323.4516

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 19999
$ws.Range("J86").Value = 19999
$ws.Range("L86").Value = 19999
$ws.Range("N86").Value = -22371
$ws.Range("H89").Value = 19999
$ws.Range("J89").Value = 19999
$ws.Range("L89").Value = 59997
$ws.Range("N89").Value = -71853
$ws.Range("H122").Value = 121213310
$ws.Range("I122").Value = 47620250
$ws.Range("J122").Value = 250001170
$ws.Range("K122").Value = 142860750
$ws.Range("L122").Value = 750003510
$ws.Range("M122").Value = -142858300
$ws.Range("N122").Value = -750008410
$ws.Range("H132").Value = 50324.094
$ws.Range("I132").Value = 41997.04
$ws.Range("J132").Value = 86408
$ws.Range("K132").Value = 125991.12
$ws.Range("L132").Value = 259224
$ws.Range("M132").Value = -123461.12
$ws.Range("N132").Value = -264284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5857.2856
$ws.Range("I7").Value = 6080.875
$ws.Range("K7").Value = 6080.875
$ws.Range("M7").Value = -5968.875
$ws.Range("H40").Value = 3047.3823
$ws.Range("I40").Value = 1946.8334
$ws.Range("J40").Value = 3647.682
$ws.Range("K40").Value = 1946.8334
$ws.Range("L40").Value = 3647.682
$ws.Range("M40").Value = -1810.8334
$ws.Range("N40").Value = -3919.682
$ws.Range("H46").Value = 1087.8889
$ws.Range("I46").Value = 718.2
$ws.Range("J46").Value = 1550
$ws.Range("K46").Value = 718.2
$ws.Range("L46").Value = 1550
$ws.Range("M46").Value = -530.2
$ws.Range("N46").Value = -1926
$ws.Range("H126").Value = 5857.2856
$ws.Range("I126").Value = 6080.875
$ws.Range("K126").Value = 18242.625
$ws.Range("M126").Value = -15772.625
$ws.Range("H132").Value = 2108.2903
$ws.Range("I132").Value = 1515.6086
$ws.Range("K132").Value = 4546.825800000001
$ws.Range("M132").Value = -2016.825800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1235.0714
$ws.Range("I126").Value = 1334.0588
$ws.Range("K126").Value = 4002.1764
$ws.Range("M126").Value = -1532.1764
$ws.Range("I136").Value = 37038652
$ws.Range("J136").Value = 1947.0625
$ws.Range("K136").Value = 111115956
$ws.Range("L136").Value = 5841.1875
$ws.Range("M136").Value = -111113406
$ws.Range("N136").Value = -10941.1875

Write-Host "Done applying changes"
